$d = $word.ActiveDocument

# 1) Rework the "Access should be enabled ..." sentence text.
$d.Content.Find.Execute(
    "Access should be enabled for both index.html and error.html pages.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Access should be enabled for index.html page.", 2) | Out-Null

# 2) Re-join the "Enable encryption ..." sentence into a single run
#    (it currently spans the old bookmark split).
$d.Content.Find.Execute(
    "Enable encryption on both files before uploading.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enable encryption on both files before uploading.", 2) | Out-Null

# Locate the "Access should be enabled ..." paragraph (and the empty,
# numId=0/tabs paragraph right after it) by content instead of a hard-coded
# index, so the script is resilient to any paragraph re-numbering.
$accessPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Access should be enabled for index.html page*") {
        $accessPara = $para
        $accessIndex = $i
        break
    }
}

# 3) Move the "_GoBack" bookmark from the "Enable encryption" paragraph to the
#    "Access should be enabled" paragraph, splitting it right before the
#    trailing period (Word automatically relocates a bookmark added again
#    under the same name).
$accessRange = $accessPara.Range
$splitPos = $accessRange.Start + 44
$insertPoint = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null

# 4) Restore the explicit <w:ilvl w:val="0"/> on the empty paragraph that
#    follows (numId stays 0 / "no list").
$tabsPara = $d.Paragraphs.Item($accessIndex + 1)
$tabsPara.Range.ListFormat.RemoveNumbers() | Out-Null
